$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.698.98"
$ws.Range("E2").Value = "  +1.00%  "

# Row 3
$ws.Range("D3").Value = "1.808.61"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.38"
$ws.Range("E5").Value = "  -0.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("E6").Value = "  +1.14%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.52"
$ws.Range("E8").Value = "  +11.60%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  -0.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0673"
$ws.Range("E10").Value = "  -2.27%  "

# Row 11
$ws.Range("E11").Value = "  +4.08%  "

# Row 12
$ws.Range("D12").Value = "2.071.17"
$ws.Range("E12").Value = "  +0.56%  "

# Row 13
$ws.Range("D13").Value = "1.803.39"
$ws.Range("E13").Value = "  +0.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.93"
$ws.Range("E14").Value = "  -3.53%  "

# Row 15
$ws.Range("E15").Value = "  -0.69%  "

# Row 16
$ws.Range("D16").Value = "34.721.97"
$ws.Range("E16").Value = "  +1.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.40"
$ws.Range("E17").Value = "  -0.47%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.95"
$ws.Range("E18").Value = "  -2.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.38"
$ws.Range("E19").Value = "  -0.92%  "

# Row 20
$ws.Range("E20").Value = "  -1.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("E21").Value = "  -3.04%  "

# Row 22
$ws.Range("E22").Value = "  +0.16%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  -1.17%  "

# Row 24
$ws.Range("E24").Value = "  -3.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.89"
$ws.Range("E25").Value = "  +0.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.72"
$ws.Range("E26").Value = "  -4.80%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.51"
$ws.Range("E27").Value = "  +0.91%  "

# Row 28
$ws.Range("E28").Value = "  +0.21%  "

# Row 29
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.78"
$ws.Range("E30").Value = "  -0.93%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.22"
$ws.Range("E31").Value = "  -1.09%  "

# Row 32
$ws.Range("E32").Value = "  -0.65%  "

# Row 33
$ws.Range("E33").Value = "  -2.31%  "

# Row 34
$ws.Range("E34").Value = "  +0.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.643"
$ws.Range("E35").Value = "  -1.91%  "

# Row 36
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.309.57"
$ws.Range("E36").Value = "  -4.75%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  +0.22%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.98"
$ws.Range("E38").Value = "  +11.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0189"
$ws.Range("E39").Value = "  +0.97%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.35"
$ws.Range("E40").Value = "  -0.48%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "84.58"
$ws.Range("E41").Value = "  +3.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.25"
$ws.Range("E42").Value = "  +6.14%  "

# Row 43
$ws.Range("E43").Value = "  +0.63%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  +0.51%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.945"
$ws.Range("E45").Value = "  -0.06%  "

# Row 46
$ws.Range("E46").Value = "  +5.11%  "

# Row 47
$ws.Range("D47").Value = "1.969.57"
$ws.Range("E47").Value = "  +0.43%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.74"
$ws.Range("E48").Value = "  -2.21%  "

# Row 49
$ws.Range("E49").Value = "  +0.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.24"
$ws.Range("E50").Value = "  -1.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0613"
$ws.Range("E51").Value = "  +0.76%  "
